# Auto-generated edit script applying crypto price/volume updates
# (commit: Updated cryptos list on Wed Sep 27 09:35:02 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.281.70"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.595.45"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.819.24"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "1.595.15"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "26.275.67"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "216.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "1.433.19"
$ws.Range("E33").Value = "  +7.25%  "
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.557"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.759"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "1.732.57"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.907"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.52%  "
